$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B6").Value = "B-"
$ws.Range("B7").Select()
